$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(43, 8).Value = 794.3333
$ws.Cells.Item(43, 9).Value = 716.6667
$ws.Cells.Item(43, 10).Value = 833.1667
$ws.Cells.Item(43, 11).Value = 716.6667
$ws.Cells.Item(43, 12).Value = 833.1667
$ws.Cells.Item(43, 13).Value = -647.6667
$ws.Cells.Item(43, 14).Value = -971.1667

$ws.Cells.Item(121, 8).Value = 6375.619
$ws.Cells.Item(121, 10).Value = 6375.619
$ws.Cells.Item(121, 12).Value = 19126.857
$ws.Cells.Item(121, 14).Value = -22620.857

$ws.Cells.Item(129, 8).Value = 841.2143
$ws.Cells.Item(129, 9).Value = 499.5
$ws.Cells.Item(129, 10).Value = 898.1667
$ws.Cells.Item(129, 11).Value = 1498.5
$ws.Cells.Item(129, 12).Value = 2694.5001
$ws.Cells.Item(129, 13).Value = 3501.5
$ws.Cells.Item(129, 14).Value = -12694.5001

$ws.Cells.Item(132, 8).Value = 43257.96
$ws.Cells.Item(132, 9).Value = 44852.043
$ws.Cells.Item(132, 10).Value = 5000
$ws.Cells.Item(132, 11).Value = 134556.129
$ws.Cells.Item(132, 12).Value = 15000
$ws.Cells.Item(132, 13).Value = -132026.129
$ws.Cells.Item(132, 14).Value = -20060

$ws.Cells.Item(135, 8).Value = 35725950
$ws.Cells.Item(135, 9).Value = 1594.5555
$ws.Cells.Item(135, 10).Value = 100029800
$ws.Cells.Item(135, 11).Value = 14350.9995
$ws.Cells.Item(135, 12).Value = 900268200
$ws.Cells.Item(135, 13).Value = -11815.9995
$ws.Cells.Item(135, 14).Value = -900273270

$ws.Cells.Item(136, 8).Value = 49866.668
$ws.Cells.Item(136, 10).Value = 49866.668
$ws.Cells.Item(136, 12).Value = 49866.668
$ws.Cells.Item(136, 14).Value = -60066.668

$ws.Cells.Item(141, 8).Value = 1947
$ws.Cells.Item(141, 9).Value = 1386.5385
$ws.Cells.Item(141, 10).Value = 4375.6665
$ws.Cells.Item(141, 11).Value = 4159.6155
$ws.Cells.Item(141, 12).Value = 13126.9995
$ws.Cells.Item(141, 13).Value = 1020.3845
$ws.Cells.Item(141, 14).Value = -23486.9995

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(25, 8).Value = 440
$ws.Cells.Item(25, 9).Value = 440
$ws.Cells.Item(25, 11).Value = 440
$ws.Cells.Item(25, 13).Value = -38

$ws.Cells.Item(61, 8).Value = 403093.12
$ws.Cells.Item(61, 9).Value = 602389.2
$ws.Cells.Item(61, 11).Value = 602389.2
$ws.Cells.Item(61, 13).Value = -602177.2

$ws.Cells.Item(102, 8).Value = 4887.5
$ws.Cells.Item(102, 9).Value = 3617.8
$ws.Cells.Item(102, 10).Value = 7003.6665
$ws.Cells.Item(102, 11).Value = 3617.8
$ws.Cells.Item(102, 12).Value = 7003.6665
$ws.Cells.Item(102, 13).Value = -1995.8
$ws.Cells.Item(102, 14).Value = -10247.6665

$ws.Cells.Item(132, 8).Value = 44243.668
$ws.Cells.Item(132, 9).Value = 2811.2727
$ws.Cells.Item(132, 11).Value = 8433.8181
$ws.Cells.Item(132, 13).Value = -5903.8181

$ws.Cells.Item(136, 8).Value = 403093.12
$ws.Cells.Item(136, 9).Value = 602389.2
$ws.Cells.Item(136, 11).Value = 1807167.6
$ws.Cells.Item(136, 13).Value = -1804617.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 2377.25
$ws.Cells.Item(20, 9).Value = 2748.923
$ws.Cells.Item(20, 11).Value = 2748.923
$ws.Cells.Item(20, 13).Value = -2501.923

$ws.Cells.Item(134, 8).Value = 33721.344
$ws.Cells.Item(134, 9).Value = 35892.766
$ws.Cells.Item(134, 10).Value = 1150
$ws.Cells.Item(134, 11).Value = 107678.298
$ws.Cells.Item(134, 12).Value = 3450
$ws.Cells.Item(134, 13).Value = -105143.298
$ws.Cells.Item(134, 14).Value = -8520

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 11782
$ws.Cells.Item(31, 10).Value = 2427.15
$ws.Cells.Item(31, 12).Value = 2427.15
$ws.Cells.Item(31, 14).Value = -3017.15

$ws.Cells.Item(34, 8).Value = 11782
$ws.Cells.Item(34, 10).Value = 2427.15
$ws.Cells.Item(34, 12).Value = 2427.15
$ws.Cells.Item(34, 14).Value = -2831.15

$ws.Cells.Item(134, 8).Value = 7660.6
$ws.Cells.Item(134, 9).Value = 941.5
$ws.Cells.Item(134, 10).Value = 21098.8
$ws.Cells.Item(134, 11).Value = 2824.5
$ws.Cells.Item(134, 12).Value = 63296.39999999999
$ws.Cells.Item(134, 13).Value = -289.5
$ws.Cells.Item(134, 14).Value = -68366.39999999999

$ws.Cells.Item(135, 8).Value = 52426.668
$ws.Cells.Item(135, 10).Value = 52426.668
$ws.Cells.Item(135, 12).Value = 52426.668
$ws.Cells.Item(135, 14).Value = -62566.668

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 7713.385
$ws.Cells.Item(2, 9).Value = 12523.75
$ws.Cells.Item(2, 10).Value = 16.8
$ws.Cells.Item(2, 11).Value = 75142.5
$ws.Cells.Item(2, 12).Value = 100.8
$ws.Cells.Item(2, 13).Value = -75029.5
$ws.Cells.Item(2, 14).Value = -326.8

$ws.Cells.Item(3, 8).Value = 1466.6666
$ws.Cells.Item(3, 9).Value = 1466.6666
$ws.Cells.Item(3, 11).Value = 4399.9998
$ws.Cells.Item(3, 13).Value = -4287.9998

$ws.Cells.Item(5, 8).Value = 813.56757
$ws.Cells.Item(5, 9).Value = 741.6923
$ws.Cells.Item(5, 10).Value = 852.5
$ws.Cells.Item(5, 11).Value = 2225.0769
$ws.Cells.Item(5, 12).Value = 2557.5
$ws.Cells.Item(5, 13).Value = -2113.0769
$ws.Cells.Item(5, 14).Value = -2781.5

$ws.Cells.Item(34, 8).Value = 641.1111
$ws.Cells.Item(34, 9).Value = 156.66667
$ws.Cells.Item(34, 10).Value = 883.3333
$ws.Cells.Item(34, 11).Value = 470.00001
$ws.Cells.Item(34, 12).Value = 2649.9999
$ws.Cells.Item(34, 13).Value = -386.00001
$ws.Cells.Item(34, 14).Value = -2817.9999

$ws.Cells.Item(68, 8).Value = 3680.6667
$ws.Cells.Item(68, 9).Value = 900
$ws.Cells.Item(68, 10).Value = 3912.389
$ws.Cells.Item(68, 11).Value = 2700
$ws.Cells.Item(68, 12).Value = 11737.167
$ws.Cells.Item(68, 13).Value = -1889
$ws.Cells.Item(68, 14).Value = -13359.167

$ws.Cells.Item(71, 8).Value = 3680.6667
$ws.Cells.Item(71, 9).Value = 900
$ws.Cells.Item(71, 10).Value = 3912.389
$ws.Cells.Item(71, 11).Value = 8100
$ws.Cells.Item(71, 12).Value = 35211.501
$ws.Cells.Item(71, 13).Value = -4044
$ws.Cells.Item(71, 14).Value = -43323.501

$ws.Cells.Item(107, 8).Value = 5157.0835
$ws.Cells.Item(107, 9).Value = 33765.332
$ws.Cells.Item(107, 10).Value = 1070.1904
$ws.Cells.Item(107, 11).Value = 101295.996
$ws.Cells.Item(107, 12).Value = 3210.5712
$ws.Cells.Item(107, 13).Value = -99375.99600000001
$ws.Cells.Item(107, 14).Value = -7050.5712

$ws.Cells.Item(135, 8).Value = 813.56757
$ws.Cells.Item(135, 9).Value = 741.6923
$ws.Cells.Item(135, 10).Value = 852.5
$ws.Cells.Item(135, 11).Value = 6675.2307
$ws.Cells.Item(135, 12).Value = 7672.5
$ws.Cells.Item(135, 13).Value = -4140.2307
$ws.Cells.Item(135, 14).Value = -12742.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 10418688
$ws.Cells.Item(70, 9).Value = 1375
$ws.Cells.Item(70, 10).Value = 15627344
$ws.Cells.Item(70, 11).Value = 1375
$ws.Cells.Item(70, 12).Value = 15627344
$ws.Cells.Item(70, 13).Value = -1105
$ws.Cells.Item(70, 14).Value = -15627884

$ws.Cells.Item(73, 8).Value = 10418688
$ws.Cells.Item(73, 9).Value = 1375
$ws.Cells.Item(73, 10).Value = 15627344
$ws.Cells.Item(73, 11).Value = 1375
$ws.Cells.Item(73, 12).Value = 15627344
$ws.Cells.Item(73, 13).Value = -439
$ws.Cells.Item(73, 14).Value = -15629216

$ws.Cells.Item(80, 8).Value = 8452.096
$ws.Cells.Item(80, 9).Value = 15101.875
$ws.Cells.Item(80, 10).Value = 4359.923
$ws.Cells.Item(80, 11).Value = 15101.875
$ws.Cells.Item(80, 12).Value = 4359.923
$ws.Cells.Item(80, 13).Value = -14103.875
$ws.Cells.Item(80, 14).Value = -6355.923

$ws.Cells.Item(83, 8).Value = 8452.096
$ws.Cells.Item(83, 9).Value = 15101.875
$ws.Cells.Item(83, 10).Value = 4359.923
$ws.Cells.Item(83, 11).Value = 75509.375
$ws.Cells.Item(83, 12).Value = 21799.615
$ws.Cells.Item(83, 13).Value = -70517.375
$ws.Cells.Item(83, 14).Value = -31783.615

$ws.Cells.Item(122, 8).Value = 5189.6
$ws.Cells.Item(122, 9).Value = 3781.2
$ws.Cells.Item(122, 10).Value = 6598
$ws.Cells.Item(122, 11).Value = 11343.6
$ws.Cells.Item(122, 12).Value = 19794
$ws.Cells.Item(122, 13).Value = -8893.599999999999
$ws.Cells.Item(122, 14).Value = -24694

$ws.Cells.Item(126, 8).Value = 6202.48
$ws.Cells.Item(126, 9).Value = 5964.1333
$ws.Cells.Item(126, 10).Value = 6560
$ws.Cells.Item(126, 11).Value = 17892.3999
$ws.Cells.Item(126, 12).Value = 19680
$ws.Cells.Item(126, 13).Value = -15422.3999
$ws.Cells.Item(126, 14).Value = -24620

$ws.Cells.Item(132, 8).Value = 120469.305
$ws.Cells.Item(132, 9).Value = 131012.625
$ws.Cells.Item(132, 10).Value = 103600
$ws.Cells.Item(132, 11).Value = 393037.875
$ws.Cells.Item(132, 12).Value = 310800
$ws.Cells.Item(132, 13).Value = -390507.875
$ws.Cells.Item(132, 14).Value = -315860

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 2725
$ws.Cells.Item(16, 9).Value = 2725
$ws.Cells.Item(16, 11).Value = 2725
$ws.Cells.Item(16, 13).Value = -2555

$ws.Cells.Item(24, 8).Value = 3500
$ws.Cells.Item(24, 10).Value = 3500
$ws.Cells.Item(24, 12).Value = 3500
$ws.Cells.Item(24, 14).Value = -4186

$ws.Cells.Item(119, 8).Value = 30420
$ws.Cells.Item(119, 10).Value = 30420
$ws.Cells.Item(119, 12).Value = 30420
$ws.Cells.Item(119, 14).Value = -40096

$ws.Cells.Item(132, 8).Value = 1608.2094
$ws.Cells.Item(132, 9).Value = 1212.2333
$ws.Cells.Item(132, 10).Value = 2522
$ws.Cells.Item(132, 11).Value = 3636.699900000001
$ws.Cells.Item(132, 12).Value = 7566
$ws.Cells.Item(132, 13).Value = -1106.699900000001
$ws.Cells.Item(132, 14).Value = -12626

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(46, 8).Value = 30000
$ws.Cells.Item(46, 10).Value = 30000
$ws.Cells.Item(46, 12).Value = 30000
$ws.Cells.Item(46, 14).Value = -30462

$ws.Cells.Item(81, 8).Value = 1952.5
$ws.Cells.Item(81, 9).Value = 2136.6667
$ws.Cells.Item(81, 10).Value = 1400
$ws.Cells.Item(81, 11).Value = 4273.3334
$ws.Cells.Item(81, 12).Value = 2800
$ws.Cells.Item(81, 13).Value = -3212.3334
$ws.Cells.Item(81, 14).Value = -4922

$ws.Cells.Item(84, 8).Value = 1952.5
$ws.Cells.Item(84, 9).Value = 2136.6667
$ws.Cells.Item(84, 10).Value = 1400
$ws.Cells.Item(84, 11).Value = 21366.667
$ws.Cells.Item(84, 12).Value = 14000
$ws.Cells.Item(84, 13).Value = -16062.667
$ws.Cells.Item(84, 14).Value = -24608

$ws.Cells.Item(100, 8).Value = 1248.5834
$ws.Cells.Item(100, 9).Value = 700
$ws.Cells.Item(100, 11).Value = 1400
$ws.Cells.Item(100, 13).Value = -859

$ws.Cells.Item(132, 8).Value = 1993.2963
$ws.Cells.Item(132, 9).Value = 1796.05
$ws.Cells.Item(132, 10).Value = 2556.8572
$ws.Cells.Item(132, 11).Value = 5388.15
$ws.Cells.Item(132, 12).Value = 7670.571599999999
$ws.Cells.Item(132, 13).Value = -2858.15
$ws.Cells.Item(132, 14).Value = -12730.5716

$ws.Cells.Item(134, 8).Value = 30000
$ws.Cells.Item(134, 10).Value = 30000
$ws.Cells.Item(134, 12).Value = 90000
$ws.Cells.Item(134, 14).Value = -95070

$ws.Cells.Item(136, 8).Value = 1439.8846
$ws.Cells.Item(136, 9).Value = 852.75
$ws.Cells.Item(136, 10).Value = 1943.1428
$ws.Cells.Item(136, 11).Value = 2558.25
$ws.Cells.Item(136, 12).Value = 5829.428400000001
$ws.Cells.Item(136, 13).Value = -8.25
$ws.Cells.Item(136, 14).Value = -10929.4284
